$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.456.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "'1.646.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").Value = "'299.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.62%  "
$ws.Range("D7").Value = "'0.3793"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").Value = "'0.3565"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").Value = "'50.11"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.49%  "
$ws.Range("D10").Value = "'0.08078"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("D11").Value = "'1.216"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.63%  "
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").Value = "'21.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.65%  "
$ws.Range("D14").Value = "'6.385"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.31%  "
$ws.Range("D15").Value = "'7.369"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").Value = "'0.00001192"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("D17").Value = "'1.644.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").Value = "'96.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").Value = "'0.06972"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "'6.766"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").Value = "'17.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.66%  "
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").Value = "'12.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("D24").Value = "'23.480.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.88%  "
$ws.Range("D25").Value = "'2.488"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("D26").Value = "'2.885"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.13%  "
$ws.Range("D27").Value = "'20.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("D28").Value = "'152.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("D29").Value = "'5.206"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("D30").Value = "'133.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("D31").Value = "'1.826.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").Value = "'6.888"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("D33").Value = "'2.137"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.93%  "
$ws.Range("D34").Value = "'11.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.69%  "
$ws.Range("D35").Value = "'1.019"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.75%  "
$ws.Range("D36").Value = "'0.02723"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.36%  "
$ws.Range("D37").Value = "'0.08723"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.18%  "
$ws.Range("D38").Value = "'0.2432"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.56%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'5.922"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "'13.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.50%  "
$ws.Range("D41").Value = "'0.06782"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.80%  "
$ws.Range("D42").Value = "'0.6875"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.91%  "
$ws.Range("D43").Value = "'1.317"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("D44").Value = "'15.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.80%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "'1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.6396"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("D47").Value = "'2.255"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.68%  "
$ws.Range("D48").Value = "'3.921"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("D49").Value = "'0.07750"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.91%  "
$ws.Range("D50").Value = "'127.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("D51").Value = "'1.153"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.58%  "
